$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" footer field from 9/16/24
#    to 10/31/24 across the slide master and every slide layout.
# ------------------------------------------------------------------
$oldDate = "9/16/24"
$newDate = "10/31/24"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ------------------------------------------------------------------
# 2) Correct the swapped plot-6 "Mai" location labels on slide 7:
#    the textbox that read TT24_231 should read TT24_232 and
#    vice versa.
# ------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
for ($i = 1; $i -le $slide7.Shapes.Count; $i++) {
    $shp = $slide7.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $t = $shp.TextFrame.TextRange.Text
        if ($t -eq "TT24_231") {
            $shp.TextFrame.TextRange.Text = "TT24_232"
        } elseif ($t -eq "TT24_232") {
            $shp.TextFrame.TextRange.Text = "TT24_231"
        }
    }
}
